# Add 26 more KTANE modules to the modules-config-details sheet (rows 156-181).
# Columns: A=Module Name, B=Module Id, C=Difficulty, D=PDF URL, E=Author(s), F=Date added (text), G=constant 1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("modules-config-details")

# Column F holds text-formatted dates (e.g. "2018-01-01"); force Text format
# *before* assigning values so Excel doesn't auto-convert them to date serials.
$ws.Range("F156:F181").NumberFormat = "@"

$rows = @(
    @('Polyhedral Maze',      'PolyhedralMazeModule',   5, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Polyhedral%20Maze.pdf',            'Timwi',                         '2018-01-01'),
    @('Symbolic Coordinates', 'symbolicCoordinates',    2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Symbolic%20Coordinates.pdf',       'Royal_Flu$h',                   '2018-01-09'),
    @('Poker',                'Poker',                  3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Poker.pdf',                       'Royal_Flu$h',                   '2018-01-09'),
    @('Sonic the Hedgehog',   'sonic',                  3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Sonic%20the%20Hedgehog.pdf',       'Royal_Flu$h',                   '2018-01-14'),
    @('Poetry',               'poetry',                 3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Poetry.pdf',                      'clutterArranger',               '2018-01-14'),
    @('Button Sequence',      'buttonSequencesModule',  3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Button%20Sequence.pdf',            'ZekNikZ',                        '2018-01-15'),
    @('Algebra',               'algebra',                3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Algebra.pdf',                    'Royal_Flu$h',                   '2018-01-17'),
    @('Visual Impairment',    'visual_impairment',      3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Visual%20Impairment.pdf',          'KingBranBran',                   '2018-01-20'),
    @('Jukebox',               'jukebox',                2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/The%20Jukebox.pdf',               'Royal_Flu$h',                   '2018-01-23'),
    @('Identity Parade',      'identityParade',         2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Identity%20Parade.pdf',            'Royal_Flu$h',                   '2018-01-25'),
    @('Maintenance',           'maintenance',            4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Maintenance.pdf',                 'Royal_Flu$h',                   '2018-01-30'),
    @('Blind Maze',            'BlindMaze',              4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Blind%20Maze.pdf',                'Riverbui, McNiko67',             '2018-01-30'),
    @('Backgrounds',           'Backgrounds',            2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Backgrounds.pdf',                 'McNiko67',                        '2018-01-30'),
    @('Mortal Kombat',        'mortalKombat',           3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Mortal%20Kombat.pdf',              'Royal_Flu$h',                   '2018-02-02'),
    @('Mashematics',          'mashematics',            2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Mashematics.pdf',                  'Marksam32',                       '2018-02-03'),
    @('Faulty Backgrounds',   'FaultyBackgrounds',      2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Faulty%20Backgrounds.pdf',         'McNiko67',                        '2018-02-03'),
    @('Radiator',              'radiator',               2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Radiator.pdf',                   'red031000, Inova',               '2018-02-04'),
    @('Modern Cipher',        'modernCipher',           3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Modern%20Cipher.pdf',              'TheFe',                           '2018-02-04'),
    @('LED Grid',              'LED Grid',               2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/LED%20Grid.pdf',                  'Royal_Flu$h',                   '2018-02-06'),
    @('Sink',                  'Sink',                   2, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Sink.pdf',                        'McNiko67',                        '2018-02-09'),
    @('iPhone',                'iPhone',                 4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/The%20iPhone.pdf',                'Royal_Flu$h',                   '2018-02-20'),
    @('Swan',                  'theSwan',                4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/The%20Swan.pdf',                  'Royal_Flu$h',                   '2018-02-22'),
    @('Waste Management',     'wastemanagement',        4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Waste%20Management.pdf',           'red031000, Inova, AppleSlice',   '2018-02-23'),
    @('Human Resources',      'HumanResourcesModule',   3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Human%20Resources.pdf',            'Elias, Timwi',                    '2018-02-26'),
    @('Skyrim',                'skyrim',                 4, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Skyrim.pdf',                      'Royal_Flu$h',                   '2018-03-02'),
    @('Burglar Alarm',        'burglarAlarm',           3, 'https://raw.githubusercontent.com/Timwi/KtaneContent/master/PDF/Burglar%20Alarm.pdf',              'Marksam32',                       '2018-03-03')
)

$r = 156
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = 1
    $r = $r + 1
}

# Mirror the author's final selection state.
$ws.Range("D182").Select()

Write-Output "Added $($rows.Count) module rows (156-181)."
